$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (13) to the color legend table: a swatch cell (A13) using a
# light-grey fill/font color, and a label cell (B13) reading "cell color".
$ws.Range("B13").Value = "cell color"
$ws.Rows.Item(13).RowHeight = 43

$swatch = $ws.Range("A13")
$swatch.Interior.Color = 16119285   # RGB(245,245,245) -> FFF5F5F5
$swatch.Font.Color = 16119285       # RGB(245,245,245) -> FFF5F5F5
$swatch.HorizontalAlignment = -4108 # xlCenter
$swatch.VerticalAlignment = -4108   # xlCenter

# Match the label cell styling used by the rest of the legend (column B).
$label = $ws.Range("B13")
$label.HorizontalAlignment = -4108
$label.VerticalAlignment = -4108

# Update the view so the new row is visible / selected, similar to the
# author's saved view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("D14").Select() | Out-Null
